$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.307.14"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.799.00"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "226.89"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "0.573"
$ws.Range("E6").Value = "  +3.51%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "35.89"
$ws.Range("E8").Value = "  +10.08%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("D12").Value = "2.060.37"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "11.56"
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("D14").Value = "1.827.64"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "0.640"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "4.47"
$ws.Range("E16").Value = "  +4.76%  "
$ws.Range("D17").Value = "34.300.58"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "68.84"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "244.39"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "0.0₃0792"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "11.51"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "4.14"
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("D25").Value = "171.77"
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  +8.62%  "
$ws.Range("D27").Value = "16.80"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").Value = "0.117"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "1.397.08"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "0.669"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "2.45"
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +10.03%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "82.27"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.956"
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "13.35"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("E46").Value = "  -3.38%  "
$ws.Range("D47").Value = "6.02"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").Value = "1.961.09"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "103.82"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("E51").Value = "  -0.12%  "
